$wb = $excel.ActiveWorkbook

# --- Update "Location Implementation" sheet data (fix bug in SE) ---
$ws = $wb.Worksheets.Item("Location Implementation")

$ws.Range("E3").Value = 0.5
$ws.Range("F3").Value = 0

$ws.Range("E4").Value = 0.5
$ws.Range("F4").Value = 0

$ws.Range("F5").Value = 0

$ws.Range("F6").Value = 0

$ws.Range("E7").Value = 5000
$ws.Range("F7").Value = 0

$ws.Range("F8").Value = 0

$ws.Range("F9").Value = 0

# Update the selection on this sheet and make it the active sheet/tab
$ws.Range("E4").Select()
$ws.Activate()
